$d = $word.ActiveDocument

# Locate the "Testing" row in the big table (Table 1) and update the
# hours tally for Mark from 5 to 7, moving the _GoBack bookmark so it
# ends up right after the new value (matching a normal in-place edit).

$tbl = $d.Tables(1)
$targetRow = $null
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $cellText = $tbl.Cell($r, 2).Range.Text
    if ($cellText -match "Testing") {
        $targetRow = $r
        break
    }
}

$cell = $tbl.Cell($targetRow, 7)

# Remove any existing _GoBack bookmark elsewhere in the document so Word
# doesn't end up with a stale duplicate once we recreate it here.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Replace "Mark:5" with "Mark:7" within this cell only.
$rng = $cell.Range
$rng.Find.Execute("Mark:5", $true, $false, $false, $false, $false, $true, 1, $false, "Mark:7", 2)

# Re-create the _GoBack bookmark right after "Mark:7" (mirrors Word's
# behaviour of dropping _GoBack at the most recent edit location).
$afterRng = $cell.Range
$afterRng.Find.Execute("Mark:7", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markEnd = $afterRng.End
$bmRange = $d.Range($markEnd, $markEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
